$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Seed rows 122:128 with the same formatting as rows 115:121 (a matching
#    7-row block), so the new rows inherit the existing visual style: the
#    date-text look of column A and the shared numeric style used in
#    columns B/D/E/F/G/H.
$ws.Range("A115:H121").Copy($ws.Range("A122:H128"))

# 2) Column A - date labels, kept as plain text like the rows above them.
#    "30/12/2022" and "31/12/2022" have a day part > 12 so they can't be
#    misread as dates; the January ones are ambiguous, so an apostrophe
#    forces literal text the same way typing it in Excel would.
$ws.Range("A122").Value = "30/12/2022"
$ws.Range("A123").Value = "31/12/2022"
$ws.Range("A124").Value = "'1/01/2023"
$ws.Range("A125").Value = "'2/01/2023"
$ws.Range("A126").Value = "'3/01/2023"
$ws.Range("A127").Value = "'4/01/2023"
$ws.Range("A128").Value = "'5/01/2023"

# Re-apply the plain (non quote-prefixed) number format from A122 onto the
# cells that needed the apostrophe, so they keep matching the look of the
# rest of the column.
$ws.Range("A122").Copy()
$ws.Range("A124:A128").PasteSpecial(-4122)

# 3) Columns B, D, E, F, G, H - raw counts for each new day.
$ws.Range("B122").Value = 2002
$ws.Range("D122").Value = 1426
$ws.Range("E122").Value = 639
$ws.Range("F122").Value = 31
$ws.Range("G122").Value = 23
$ws.Range("H122").Value = 5663

$ws.Range("B123").Value = 1424
$ws.Range("D123").Value = 993
$ws.Range("E123").Value = 642
$ws.Range("F123").Value = 28
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 4434

$ws.Range("B124").Value = 1226
$ws.Range("D124").Value = 927
$ws.Range("E124").Value = 669
$ws.Range("F124").Value = 29
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 3197

# Row 125 (2/01/2023) has no New/PCR/RAT cases reported yet, so B/C/D stay blank.
$ws.Range("B125").ClearContents()
$ws.Range("D125").ClearContents()
$ws.Range("E125").Value = 745
$ws.Range("F125").Value = 44
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 0

$ws.Range("B126").Value = 2172
$ws.Range("D126").Value = 1870
$ws.Range("E126").Value = 634
$ws.Range("F126").Value = 32
$ws.Range("G126").Value = 30
$ws.Range("H126").Value = 10950

$ws.Range("B127").Value = 1629
$ws.Range("D127").Value = 1252
$ws.Range("E127").Value = 638
$ws.Range("F127").Value = 32
$ws.Range("G127").Value = 30
$ws.Range("H127").Value = 4498

$ws.Range("B128").Value = 1650
$ws.Range("D128").Value = 1292
$ws.Range("E128").Value = 545
$ws.Range("F128").Value = 27
$ws.Range("G128").Value = 30
$ws.Range("H128").Value = 5470

# 4) Column C - "PCR cases" = New Cases minus RAT cases, the same formula
#    used throughout the sheet. Row 125 has no New/RAT data, so it stays blank.
$ws.Range("C122:C124").Formula = "=B122-D122"
$ws.Range("C125").ClearContents()
$ws.Range("C126:C128").Formula = "=B126-D126"

# 5) Leave the selection where data entry finished.
$ws.Range("A128").Select()
